# finish wk 9 and 10 schedule and PCA hw
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slo_detail")

# ---- Row 11 (week 9, Mon Mar 19) ----
# D11 SLO text is rewritten (PCA dimension-reduction objectives)
$ws.Range("D11").Value = "Explain how PCA can be used as a dimension reduction technique`nExplain the difference between multivariate and multivariable`nConduct a PCA using both the correlation and covariance matrix`n"
$ws.Rows.Item(11).RowHeight = 105

# ---- Row 12 (week 10, Mon Mar 26) ----
$ws.Range("D12").Value = "Use visualization techniques to identify the number of PC's to retain`nExplain the difference between PCA and FA`nCreate a latent factor model, visualize and interpret results. "
$ws.Range("E12").Value = "Read ASCN 15.1-15.3, and PMA6 15.1-15.4"
$ws.Range("F12").Value = "Visualizing and interpreting PC's "
$ws.Range("G12").Value = "Introduction to Factor Analysis"
$ws.Range("H12").Value = "Factor extraction & scores"

# ---- Row 13 (week 11, Mon Apr 2) ----
$ws.Range("D13").Value = "Create a latent factor model, visualize and interpret results. `nUse latent factors as a predictor in another model"
$ws.Range("E13").Value = "Read ASCM 15.4-end, PMA6 15.5-end"
$ws.Range("F13").Value = "Factor rotation"
$ws.Range("G13").Value = "Quiz. Open work day"
$ws.Range("H13").Value = "Open work day"
$ws.Rows.Item(13).RowHeight = 60

# ---- View state: slo_detail becomes the active tab, scrolled/selected ----
$ws.Activate() | Out-Null
$ws.Range("E11").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollRow = 8
